$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F334").Value = 194914
$ws.Range("G334").Value = 3415
$ws.Range("F335").Value = 130479
$ws.Range("G335").Value = 2968
$ws.Range("F336").Value = 101976
$ws.Range("F337").Value = 104232
$ws.Range("G337").Value = 2982
$ws.Range("F338").Value = 226342
$ws.Range("G338").Value = 3188
$ws.Range("F339").Value = 653498
$ws.Range("F340").Value = 380330
$ws.Range("G340").Value = 3261
$ws.Range("F341").Value = 294676
$ws.Range("G341").Value = 3662
$ws.Range("F342").Value = 178597
$ws.Range("G342").Value = 3070
$ws.Range("F343").Value = 132431
$ws.Range("G343").Value = 2972
$ws.Range("F344").Value = 135978
$ws.Range("F346").Value = 667565
$ws.Range("G346").Value = 4769
$ws.Range("F347").Value = 341307
$ws.Range("G347").Value = 2888
$ws.Range("F348").Value = 231712
$ws.Range("F349").Value = 159655
$ws.Range("G349").Value = 2747
$ws.Range("F350").Value = 127578
$ws.Range("F351").Value = 150222
$ws.Range("F352").Value = 306664
$ws.Range("G352").Value = 3548
$ws.Range("F353").Value = 718523
$ws.Range("G353").Value = 5248
$ws.Range("F354").Value = 305349
$ws.Range("G354").Value = 2783
$ws.Range("F355").Value = 221701
$ws.Range("G355").Value = 3440
$ws.Range("F356").Value = 160400
$ws.Range("F357").Value = 138317
$ws.Range("F358").Value = 157564
$ws.Range("F359").Value = 320283
$ws.Range("G359").Value = 3354
$ws.Range("F360").Value = 738891
$ws.Range("G360").Value = 5050
$ws.Range("F361").Value = 329455
$ws.Range("G361").Value = 2582
$ws.Range("F362").Value = 225462
$ws.Range("G362").Value = 3115
$ws.Range("F363").Value = 185978
$ws.Range("G363").Value = 2722
$ws.Range("F364").Value = 164710
$ws.Range("G364").Value = 2423
$ws.Range("F365").Value = 178741
$ws.Range("G365").Value = 2379
$ws.Range("F366").Value = 331141
$ws.Range("G366").Value = 2795
$ws.Range("F367").Value = 738428
$ws.Range("G367").Value = 3739
$ws.Range("F368").Value = 338095
$ws.Range("G368").Value = 2233
$ws.Range("F369").Value = 226362
$ws.Range("G369").Value = 2483
$ws.Range("F370").Value = 173224
$ws.Range("G370").Value = 1958
